# Fix truth table for LUI and AUIPC
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Truth Table")

$targets = @{
    "M1" = "Inst or rs2"
    "N1" = "ExtendSel"
    "K2" = "01000"
    "M2" = "0"
    "N2" = "10"
    "K3" = "01000"
    "M3" = "0"
    "N3" = "10"
    "M4" = "0"
    "N4" = "11"
    "M5" = "0"
    "N5" = "00"
    "K6" = "00000"
    "M6" = "1"
    "N6" = "xx"
    "T6" = "0"
    "K7" = "00000"
    "M7" = "1"
    "N7" = "xx"
    "T7" = "0"
    "K8" = "00000"
    "M8" = "1"
    "N8" = "xx"
    "T8" = "0"
    "K9" = "00000"
    "M9" = "1"
    "N9" = "xx"
    "T9" = "0"
    "K10" = "00000"
    "M10" = "1"
    "N10" = "xx"
    "T10" = "0"
    "K11" = "00000"
    "M11" = "1"
    "N11" = "xx"
    "T11" = "0"
    "M12" = "0"
    "N12" = "00"
    "M13" = "0"
    "N13" = "00"
    "M14" = "0"
    "N14" = "00"
    "M15" = "0"
    "N15" = "00"
    "M16" = "0"
    "N16" = "00"
    "M17" = "0"
    "N17" = "01"
    "M18" = "0"
    "N18" = "01"
    "M19" = "0"
    "N19" = "01"
    "M20" = "0"
    "N20" = "00"
    "M21" = "0"
    "N21" = "00"
    "M22" = "0"
    "N22" = "00"
    "M23" = "0"
    "N23" = "00"
    "M24" = "0"
    "N24" = "00"
    "M25" = "0"
    "N25" = "00"
    "M26" = "0"
    "N26" = "00"
    "M27" = "0"
    "N27" = "00"
    "M28" = "0"
    "N28" = "00"
    "M29" = "1"
    "N29" = "xx"
    "M30" = "1"
    "N30" = "xx"
    "M31" = "1"
    "N31" = "xx"
    "M32" = "1"
    "N32" = "xx"
    "M33" = "1"
    "N33" = "xx"
    "M34" = "1"
    "N34" = "xx"
    "M35" = "1"
    "N35" = "xx"
    "M36" = "1"
    "N36" = "xx"
    "M37" = "1"
    "N37" = "xx"
    "M38" = "1"
    "N38" = "xx"
}

foreach ($key in $targets.Keys) {
    $ws.Range($key).Value = $targets[$key]
}
